# Adds the 10 new "liste référence" rows (408-417) introduced by the
# upstream commit: five new heat-pump/district-heating categories
# (LE-LP, ME-LP, HE-LP, HE-LP-EBREP, HE-LP-EBREP-HYB) followed by the
# matching "-HP" variants, each reusing one of the five colour swatches
# already present on the sheet (columns A, B and C all repeat the same
# label; column D carries the colour code together with that colour's
# existing fill style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New category labels (columns A/B/C all repeat the same value) paired
# with the colour-code text that belongs in column D and the row that
# already uses the matching D-column fill style (used below to copy the
# formatting across, since the fills already exist in the workbook).
$newRows = @(
  @{ Row = 408; Label = "LE-LP";             Color = "#ffe36d"; StyleSource = 121 },
  @{ Row = 409; Label = "ME-LP";             Color = "#f6b4a4"; StyleSource = 17  },
  @{ Row = 410; Label = "HE-LP";             Color = "#9bd4dc"; StyleSource = 89  },
  @{ Row = 411; Label = "HE-LP-EBREP";       Color = "#314deb"; StyleSource = 32  },
  @{ Row = 412; Label = "HE-LP-EBREP-HYB";   Color = "#9bdb9a"; StyleSource = 275 },
  @{ Row = 413; Label = "LE-HP";             Color = "#ffe36d"; StyleSource = 121 },
  @{ Row = 414; Label = "ME-HP";             Color = "#f6b4a4"; StyleSource = 17  },
  @{ Row = 415; Label = "HE-HP";             Color = "#9bd4dc"; StyleSource = 89  },
  @{ Row = 416; Label = "HE-HP-EBREP";       Color = "#314deb"; StyleSource = 32  },
  @{ Row = 417; Label = "HE-HP-EBREP-HYB";   Color = "#9bdb9a"; StyleSource = 275 }
)

foreach ($item in $newRows) {
  $r = $item.Row

  $ws.Cells.Item($r, 1).Value = $item.Label
  $ws.Cells.Item($r, 2).Value = $item.Label
  $ws.Cells.Item($r, 3).Value = $item.Label
  $ws.Cells.Item($r, 4).Value = $item.Color

  # Copy the column-D fill/format from the existing row that already
  # uses this colour so the new cell gets the identical style.
  $ws.Cells.Item($item.StyleSource, 4).Copy()
  $ws.Cells.Item($r, 4).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Keep the sheet's active selection in sync with where Excel would have
# left it after appending these rows.
$ws.Range("D413").Select()
